# Major refactor; improving PCA
# Re-run of the OCSVM metrics: update the reported scores and the
# kernel hyper-parameters (gamma, kernel) that produced them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore/ensure the workbook's legacy indexed colour palette is the
# standard Excel 64-colour table (no-op visually, but present in the
# saved styles part of a normally-resaved workbook).
$paletteColors = @(
    0x000000,0xFFFFFF,0xFF0000,0x00FF00,0x0000FF,0xFFFF00,0xFF00FF,0x00FFFF,
    0x000000,0xFFFFFF,0xFF0000,0x00FF00,0x0000FF,0xFFFF00,0xFF00FF,0x00FFFF,
    0x800000,0x008000,0x000080,0x808000,0x800080,0x008080,0xC0C0C0,0x808080,
    0x9999FF,0x993366,0xFFFFCC,0xCCFFFF,0x660066,0xFF8080,0x0066CC,0xCCCCFF,
    0x000080,0xFF00FF,0xFFFF00,0x00FFFF,0x800080,0x800000,0x008080,0x0000FF,
    0x00CCFF,0xCCFFFF,0xCCFFCC,0xFFFF99,0x99CCFF,0xFF99CC,0xCC99FF,0xFFCC99,
    0x3366FF,0x33CCCC,0x99CC00,0xFFCC00,0xFF9900,0xFF6600,0x666699,0x969696,
    0x003366,0x339966,0x003300,0x333300,0x993300,0x993366,0x333399,0x333333
)
try {
    $wb.Colors = $paletteColors
} catch {
    # Older/limited hosts may not expose a settable palette - safe to ignore.
}

# metrics (rows map A2:A9 -> B2:B9)
$ws.Range("B2").Value = 4.473053169250488     # fit_time
$ws.Range("B3").Value = 0.009083986282348633  # score_time
$ws.Range("B4").Value = 0.7700000000000001    # accuracy
$ws.Range("B5").Value = 0.9700534759358289    # precision
$ws.Range("B6").Value = 0.5585714285714285    # recall
$ws.Range("B7").Value = 0.7700000000000001    # f1_micro
$ws.Range("B8").Value = 0.7552539755429658    # f1_macro
$ws.Range("B9").Value = 0.7558471381924038    # f1_weighted

# hyper-parameters: gamma was the string "scale", now an explicit numeric value;
# kernel switches from "rbf" to "sigmoid".
$ws.Range("B11").Value = 0.001
$ws.Range("B12").Value = "sigmoid"
